$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "34.557.40"

$ws.Range("D3").Value = "1.812.39"
$ws.Range("E3").Value = "  +0.75%  "

$ws.Range("E4").Value = "  +0.07%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "229.09"
$ws.Range("E5").Value = "  +0.65%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.581"
$ws.Range("E6").Value = "  +4.49%  "

$ws.Range("E7").Value = "  +0.10%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "34.90"
$ws.Range("E8").Value = "  +6.73%  "

$ws.Range("E9").Value = "  +1.77%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0697"
$ws.Range("E10").Value = "  +0.05%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0955"
$ws.Range("E11").Value = "  +0.48%  "

$ws.Range("D12").Value = "2.075.93"
$ws.Range("E12").Value = "  +0.86%  "

$ws.Range("D13").Value = "1.822.60"
$ws.Range("E13").Value = "  +1.44%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "11.22"
$ws.Range("E14").Value = "  +0.27%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.648"
$ws.Range("E15").Value = "  +1.49%  "

$ws.Range("E16").Value = "  +3.33%  "

$ws.Range("D17").Value = "34.543.81"
$ws.Range("E17").Value = "  -0.01%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.30"
$ws.Range("E18").Value = "  +0.62%  "

$ws.Range("E19").Value = "  -0.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "245.90"
$ws.Range("E20").Value = "  -0.52%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.48"
$ws.Range("E21").Value = "  +1.65%  "

$ws.Range("E22").Value = "  +0.07%  "

$ws.Range("E23").Value = "  -0.44%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "173.16"
$ws.Range("E24").Value = "  +2.39%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.12"
$ws.Range("E25").Value = "  +3.10%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.97"
$ws.Range("E26").Value = "  +9.45%  "

$ws.Range("E27").Value = "  +1.70%  "

$ws.Range("E28").Value = "  +3.24%  "

$ws.Range("E29").Value = "  -0.11%  "

$ws.Range("E30").Value = "  -2.42%  "

$ws.Range("E31").Value = "  +1.41%  "

$ws.Range("E32").Value = "  +1.12%  "

$ws.Range("E33").Value = "  +0.17%  "

$ws.Range("E34").Value = "  +0.52%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.684"
$ws.Range("E35").Value = "  +1.61%  "

$ws.Range("D36").Value = "1.398.77"
$ws.Range("E36").Value = "  -2.26%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.47"
$ws.Range("E37").Value = "  -4.83%  "

$ws.Range("E39").Value = "  -0.33%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "83.77"
$ws.Range("E40").Value = "  -1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.962"
$ws.Range("E41").Value = "  +2.41%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.83"
$ws.Range("E42").Value = "  +2.55%  "

$ws.Range("E43").Value = "  -0.82%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.42"
$ws.Range("E44").Value = "  -3.21%  "

$ws.Range("E45").Value = "  +4.00%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0515"
$ws.Range("E46").Value = "  -1.72%  "

$ws.Range("E47").Value = "  -1.45%  "

$ws.Range("D48").Value = "1.974.96"
$ws.Range("E48").Value = "  +0.91%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "105.30"
$ws.Range("E49").Value = "  -0.49%  "

$ws.Range("E50").Value = "  +1.79%  "

$ws.Range("E51").Value = "  +0.04%  "
